# Upload pictures for the journal
# Replace the per-row "Done" markers (col M) + inline legend text (col N)
# with a cleaned-up table (table shrinks to the real 12 data rows) and move
# the legend text down into column A below the table, trimming the leftover
# numbered filler rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- legend text moves from N1:N3 into A15:A17, N column cleared ---
$legendTop   = $ws.Range("N1").Value()
$legendIncl  = $ws.Range("N2").Value()
$legendFaked = $ws.Range("N3").Value()

$ws.Range("N1").Clear()
$ws.Range("N2").Clear()
$ws.Range("N3").Clear()

# --- "Done" markers in column M (rows 2-12) are cleared, and the same
#     (now blank) style block grows to cover M13 as well ---
$ws.Range("M2:M13").ClearContents()

# the fill used by that style changes from the old accent green (theme 9)
# to a plain background (theme 0 / white)
$ws.Range("M2:M13").Interior.ThemeColor = 2        # xlThemeColorLight1 -> theme="0"

# --- trim the leftover "Step#" filler rows 14-19 down to just the legend ---
# old A14 (lone "13") goes away completely
$ws.Range("A14").Clear()

# old A15..A17 held plain numbers (14,15,16); they now hold the legend text
# that used to live in N1:N3, and N15's helper value (3) stays put
$ws.Range("A15").Value = $legendTop
$ws.Range("A16").Value = $legendIncl
$ws.Range("A17").Value = $legendFaked

# old A18/A19 become blank, but keep their style
$ws.Range("A18").ClearContents()
$ws.Range("A19").ClearContents()

# --- table shrinks to the real data range (header + 12 rows) ---
$ws.ListObjects("Tabel5").Resize($ws.Range("A1:L13"))

# --- view/selection cosmetics ---
$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("L21").Select()
